$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the existing header cell (B1) onto the new
# header cell F1, then set its text.
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "VIMMP_DEF"

# Fill in the new VIMMP_DEF column for every data row with the default
# "empty list" placeholder value used throughout the mapping sheet.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
